$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the two new "quote type" worksheets at their correct tab
#    positions, matching the sheet order implied by the A-number in the name.
# ---------------------------------------------------------------------------

# "A05 vie sainte dieudonnee" belongs right before "A06 saint jean evangeliste"
$before05 = $wb.Worksheets.Item("A06 saint jean evangeliste")
$newA05 = $wb.Worksheets.Add($before05)
$newA05.Name = "A05 vie sainte dieudonnee"

# "A10 poines enfer" belongs right before "A11 vie saint sebastien"
$before10 = $wb.Worksheets.Item("A11 vie saint sebastien")
$newA10 = $wb.Worksheets.Add($before10)
$newA10.Name = "A10 poines enfer"

# ---------------------------------------------------------------------------
# 2. Give both brand-new sheets the same header row used by every other
#    sheet in the workbook (line_n / prev_line / line / next_line), styled
#    to match (bold, centred).
# ---------------------------------------------------------------------------
foreach ($newSheet in @($newA05, $newA10)) {
    $newSheet.Range("A1").Value = "line_n"
    $newSheet.Range("B1").Value = "prev_line"
    $newSheet.Range("C1").Value = "line"
    $newSheet.Range("D1").Value = "next_line"
    $headerRange = $newSheet.Range("A1:D1")
    $headerRange.Font.Bold = $true
    $headerRange.HorizontalAlignment = -4108
}

# ---------------------------------------------------------------------------
# 3. Fix the transcription of the saint Thibault quote (punctuation /
#    capitalisation corrections) living on the "A17 robert deable" sheet.
# ---------------------------------------------------------------------------
$wsThibault = $wb.Worksheets.Item("A17 robert deable")
$wsThibault.Range("B2").Value = [char]0x201C + "Seneschal," + [char]0x201D + " dist le roy, " + [char]0x201C + "pour le cors saint Thibaut,"
$wsThibault.Range("C2").Value = "Fustes vous la pour moi?" + [char]0x201D + " " + [char]0x201C + "Ouïl, se Diex me saut." + [char]0x201D
$wsThibault.Range("D2").Value = [char]0x201C + "Seneschal," + [char]0x201D + " dist le roy, " + [char]0x201C + "preuz estes et gentis,"

# ---------------------------------------------------------------------------
# 4. Restore the originally-active sheet/tab selection.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item(1).Activate()
